# Update countries & provincias Spain
# Applies the data refresh captured in the commit diff for paises.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last refreshed" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Junio de 2020 a las 08:25"

# --- Row 7: India ---
$ws.Range("B7").Value = 396182
$ws.Range("C7").Value = 370
$ws.Range("D7").Value = 214209
$ws.Range("E7").Value = 169003

# --- Row 76: Uzbekistan ---
$ws.Range("B76").Value = 6025
$ws.Range("C76").Value = 79
$ws.Range("E76").Value = 1733

# --- Row 86: El Salvador ---
$ws.Range("D86").Value = 2326
$ws.Range("E86").Value = 1910
$ws.Range("G86").Value = 7
$ws.Range("H86").Value = 93

# --- Row 95: Kirguistan ---
$ws.Range("B95").Value = 2981
$ws.Range("C95").Value = 192
$ws.Range("D95").Value = 1981
$ws.Range("E95").Value = 965
$ws.Range("G95").Value = 3
$ws.Range("H95").Value = 35

# --- Row 131: Jordania ---
$ws.Range("B131").Value = 898
$ws.Range("C131").Value = 2
$ws.Range("D131").Value = 752
$ws.Range("E131").Value = 132

# --- Rows 202/203: reordering of Fiyi / Dominica (totals tied at 18, order swapped) ---
$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"

# --- Rows 208/209: reordering of Santa Sede / Islas Turcas y Caicos, with updated data ---
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

$ws.Range("A209").Value = "Santa Sede"
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0
